$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New summary row: average of column J (|S*|/n) over the data rows
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# New summary block (rows 14-17): labels in column A, stats in column B
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Formatting for the new label/value cells: bold 12pt font, vertically centered,
# with a slightly taller row height (matches the taller 12pt font).
$labelRange = $ws.Range("A14:B17")
$labelRange.Font.FontStyle = "Bold"
$labelRange.Font.Size = 12
$labelRange.VerticalAlignment = -4108

$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection state as left by the editing session
$ws.Range("A14:B17").Select()

Write-Output "done"
